$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value of 45183 for rows 2-101.
# Update it to 45184 (i.e. advance the "changed" date by one day) for every
# data row that currently holds that value.
for ($row = 2; $row -le 101; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value = 45184
    }
}
